$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.046.52'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '3.523.12'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("D7").Value = '3.521.81'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.125'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.67%  '
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '4.121.48'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.70'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000182'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '3.525.36'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '64.999.15'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '392.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.581'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("B25").Value = 'WrappedeETH'
$ws.Range("C25").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D25").Value = '3.665.31'
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.34'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '3.529.50'
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.12%  '
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '168.41'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("E41").Value = '  +1.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.825'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("E43").Value = '  +6.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").Value = '2.429.29'
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.911'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.81%  '
